$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the Price column (D) stored as text -- values like "20.357.18"
# or "0.9510" are display strings, not numbers, and we do not want
# Excel to silently coerce them into doubles when assigned via .Value.
# (Row 38's price text is untouched by this update, so it is skipped.)
$ws.Range("D2:D37").NumberFormat = "@"
$ws.Range("D39:D51").NumberFormat = "@"

$data = @(
    @{Row=2; D="20.357.18"; E="  +1.93%  "},
    @{Row=3; D="1.458.62"; E="  +3.34%  "},
    @{Row=4; D="1.008"; E="  +0.59%  "},
    @{Row=5; D="0.9438"; E="  -5.69%  "},
    @{Row=6; D="274.36"; E="  -0.46%  "},
    @{Row=7; D="0.3645"; E="  -0.44%  "},
    @{Row=8; D="0.3065"; E="  -1.60%  "},
    @{Row=9; D="39.77"; E="  -0.01%  "},
    @{Row=10; D="1.031"; E="  -0.40%  "},
    @{Row=11; D="0.06549"; E="  +0.63%  "},
    @{Row=12; D="0.9984"; E="  -0.33%  "},
    @{Row=13; D="5.382"; E="  -2.57%  "},
    @{Row=14; D="17.79"; E="  +0.20%  "},
    @{Row=15; D="6.097"; E="  -1.61%  "},
    @{Row=16; D="0.00001021"; E="  +0.21%  "},
    @{Row=17; D="1.458.18"; E="  +3.25%  "},
    @{Row=18; D="0.9604"; E="  -4.01%  "},
    @{Row=19; D="0.05747"; E="  +1.01%  "},
    @{Row=20; D="69.39"; E="  -2.21%  "},
    @{Row=21; D="5.415"; E="  -3.78%  "},
    @{Row=22; D="14.39"; E="  -2.51%  "},
    @{Row=23; D="10.84"; E="  -1.56%  "},
    @{Row=24; D="2.234"; E="  -1.13%  "},
    @{Row=25; D="20.393.13"; E="  +1.97%  "},
    @{Row=26; D="140.35"; E="  +5.35%  "},
    @{Row=27; D="2.075"; E="  -8.80%  "},
    @{Row=28; D="17.06"; E="  -0.95%  "},
    @{Row=29; D="1.611.07"; E="  +2.57%  "},
    @{Row=30; D="111.50"; E="  +1.62%  "},
    @{Row=31; D="3.758"; E="  -4.93%  "},
    @{Row=32; D="4.848"; E="  -8.51%  "},
    @{Row=33; B="Stellar"; C="https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; D="0.07783"; E="  +1.20%  "},
    @{Row=34; B="ImmutableX"; C="https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D="0.7844"; E="  -4.77%  "},
    @{Row=35; D="1.500"; E="  +1.33%  "},
    @{Row=36; D="0.05706"; E="  -3.46%  "},
    @{Row=37; D="4.646"; E="  -5.85%  "},
    @{Row=38; E="  +2.74%  "},
    @{Row=39; D="0.02018"; E="  -3.13%  "},
    @{Row=40; B="Frax"; C="https://coinranking.com/coin/KfWtaeV1W+frax-frax"; D="0.9510"; E="  -4.90%  "},
    @{Row=41; B="Aptos"; C="https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; D="10.31"; E="  -1.93%  "},
    @{Row=42; D="0.1858"; E="  -2.58%  "},
    @{Row=43; D="7.389"; E="  -11.91%  "},
    @{Row=44; D="0.5245"; E="  -1.44%  "},
    @{Row=45; D="3.481"; E="  -1.50%  "},
    @{Row=46; D="11.83"; E="  -4.36%  "},
    @{Row=47; D="116.61"; E="  +1.33%  "},
    @{Row=48; D="0.5130"; E="  -1.53%  "},
    @{Row=49; D="1.742"; E="  -1.81%  "},
    @{Row=50; D="0.06403"; E="  +3.42%  "},
    @{Row=51; D="0.9880"; E="  -1.15%  "}
)

foreach ($item in $data) {
    $r = $item.Row
    if ($item.ContainsKey("B")) { $ws.Cells.Item($r, 2).Value = $item.B }
    if ($item.ContainsKey("C")) { $ws.Cells.Item($r, 3).Value = $item.C }
    if ($item.ContainsKey("D")) { $ws.Cells.Item($r, 4).Value = $item.D }
    $ws.Cells.Item($r, 5).Value = $item.E
}
